$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46053
$ws.Range("B2").Value = 0.01
$ws.Range("C2").Value = -0.01
$ws.Range("D2").Value = -0.01
$ws.Range("E2").Value = -0.02
$ws.Range("F2").Value = -0.03
$ws.Range("G2").Value = -0.42
$ws.Range("H2").Value = -0.19
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = -0.11
$ws.Range("L2").Value = -0.35
$ws.Range("M2").Value = -0.44
$ws.Range("N2").Value = -0.42
$ws.Range("O2").Value = -0.2
$ws.Range("P2").Value = -0.18
$ws.Range("Q2").Value = -0.17
$ws.Range("R2").Value = -0.04
$ws.Range("S2").Value = 0.44
$ws.Range("T2").Value = 9.039999999999999
$ws.Range("U2").Value = 29.17
$ws.Range("V2").Value = 35
$ws.Range("W2").Value = 32.21
$ws.Range("X2").Value = 21.45
$ws.Range("Y2").Value = 6.91
$ws.Range("Z2").Value = 5.48
$ws.Range("AB2").Value = 23.89
$ws.Range("AD2").Value = 33.61
$ws.Range("AF2").Value = 19.1
$ws.Range("AG2").Value = "0h-17h"
